# New crime data collected - weekly CompStat update (90th Precinct)
# Report header: bulletin number 24 -> 25, week 6/10-6/16/2024 -> 6/17-6/23/2024

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates
# ---------------------------------------------------------------------------
$ws.Range("A8").Value  = "Volume 31   Number  25"
$ws.Range("C9").Value  = "Report Covering the Week  6/17/2024  Through  6/23/2024"

# ---------------------------------------------------------------------------
# Column H got wider (bestFit recalculated after the new figures landed)
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# ---------------------------------------------------------------------------
# Helper: force a cell to the "blank dash" text style used throughout this
# sheet for zero/undefined figures ("0" or "***.*" shared strings) without
# leaving a quote-prefix behind.
# ---------------------------------------------------------------------------
function Set-DashText($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
Set-DashText "F14" "0"
$ws.Range("N14").Value = -70

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
Set-DashText "C15" "0"
$ws.Range("N15").Value = -52.380952380952

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 17
$ws.Range("H16").Value = 70
$ws.Range("I16").Value = 96
$ws.Range("J16").Value = 86
$ws.Range("K16").Value = 11.627906976744
$ws.Range("L16").Value = 18.518518518518
$ws.Range("M16").Value = -39.622641509434
$ws.Range("N16").Value = -83.275261324041

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 50
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 40.909090909090
$ws.Range("I17").Value = 148
$ws.Range("J17").Value = 141
$ws.Range("K17").Value = 4.964539007092
$ws.Range("L17").Value = 19.354838709677
$ws.Range("M17").Value = 66.292134831460
$ws.Range("N17").Value = -39.094650205761

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 26
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 62.5
$ws.Range("I18").Value = 169
$ws.Range("J18").Value = 101
$ws.Range("K18").Value = 67.326732673267
$ws.Range("L18").Value = 14.965986394557
$ws.Range("M18").Value = -6.111111111111
$ws.Range("N18").Value = -71.644295302013

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -47.368421052631
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 75
$ws.Range("H19").Value = -28
$ws.Range("I19").Value = 338
$ws.Range("J19").Value = 353
$ws.Range("K19").Value = -4.249291784702
$ws.Range("L19").Value = 18.181818181818
$ws.Range("M19").Value = 66.502463054187
$ws.Range("N19").Value = 38.524590163934

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
Set-DashText "C20" "0"
Set-DashText "D20" "0"
Set-DashText "E20" "***.*"
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -25
$ws.Range("I20").Value = 62
$ws.Range("K20").Value = -12.676056338028
$ws.Range("L20").Value = -7.462686567164
$ws.Range("M20").Value = -11.428571428571
$ws.Range("N20").Value = -84.577114427860

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -22.580645161290
$ws.Range("F21").Value = 140
$ws.Range("H21").Value = 2.941176470588
$ws.Range("I21").Value = 826
$ws.Range("J21").Value = 763
$ws.Range("K21").Value = 8.256880733944
$ws.Range("L21").Value = 16.174402250351
$ws.Range("M21").Value = 17.163120567375
$ws.Range("N21").Value = -60.478468899521

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
Set-DashText "D22" "0"
Set-DashText "E22" "***.*"
Set-DashText "F22" "0"
$ws.Range("H22").Value = -100
$ws.Range("L22").Value = -47.826086956521
$ws.Range("M22").Value = -42.857142857142

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -66.666666666666
$ws.Range("G23").Value = 19
$ws.Range("H23").Value = -36.842105263157
$ws.Range("I23").Value = 88
$ws.Range("J23").Value = 95
$ws.Range("K23").Value = -7.368421052631
$ws.Range("L23").Value = 20.547945205479
$ws.Range("M23").Value = 37.5

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -4.347826086956
$ws.Range("F24").Value = 101
$ws.Range("G24").Value = 107
$ws.Range("H24").Value = -5.607476635514
$ws.Range("I24").Value = 529
$ws.Range("J24").Value = 503
$ws.Range("K24").Value = 5.168986083499
$ws.Range("L24").Value = -9.726962457337
$ws.Range("M24").Value = -4.166666666666

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = 6.25
$ws.Range("I25").Value = 177
$ws.Range("J25").Value = 63
$ws.Range("K25").Value = 180.952380952381
$ws.Range("L25").Value = 31.111111111111

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = -46.666666666666
$ws.Range("F26").Value = 50
$ws.Range("G26").Value = 51
$ws.Range("H26").Value = -1.960784313725
$ws.Range("I26").Value = 275
$ws.Range("J26").Value = 223
$ws.Range("K26").Value = 23.318385650224
$ws.Range("L26").Value = 13.168724279835
$ws.Range("M26").Value = 15.546218487395

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------------
Set-DashText "C27" "0"
$ws.Range("L27").Value = 30

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("C28").Value = 2
$ws.Range("E28").Value = 100
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 28
$ws.Range("J28").Value = 31
$ws.Range("K28").Value = -9.677419354838
$ws.Range("L28").Value = 27.272727272727

# ---------------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------------
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = 0

# ---------------------------------------------------------------------------
# Row 30 - Shooting Inc.
# ---------------------------------------------------------------------------
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = 0

# ---------------------------------------------------------------------------
# Row 31 - Hate Crimes
# ---------------------------------------------------------------------------
Set-DashText "F31" "0"
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = -100
$ws.Range("J31").Value = 6
$ws.Range("K31").Value = 150
